# Update "CIERRE GASTOS ADMINISTRATIVOS DICIEMBRE 2025.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fill RAZON_SOCIAL (column G) first for the new rows, in row order ---
$ws.Range("G23").Value = "TRANSPORTES Y MANIOBRAS S.C.R.L."
$ws.Range("G24").Value = "BUSY BIZ SOLUTIONS S.A.C."
$ws.Range("G25").Value = "FAMIP INDUSTRIAL S.A.C."
$ws.Range("G26").Value = "DORAL SOLUCIONES PERU E.I.R.L."
$ws.Range("G27").Value = "HINOSTROZA GOMEZ JOSE"
$ws.Range("G28").Value = "REPRESENTACIONES E & N PALOMINO S.A.C."
$ws.Range("G29").Value = "COPY FAST TECH S.A.C."
$ws.Range("G30").Value = "BONILLA PEREZ JHON ARMANDO"
$ws.Range("G31").Value = "MAMANI MOGROVEJO RONDOLFO"

# --- Row 23 ---
$ws.Range("A23").Value = 76474794
$ws.Range("B23").Value = "SANDRA BENAVIDES"
$ws.Range("C23").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D23").Value = "PREJUDICIAL FLUJO"
$ws.Range("E23").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F23").Value = 20505874529
$ws.Range("H23").Value = 20505874529
$ws.Range("I23").Value = 46003
$ws.Range("L23").Value = 80
$ws.Range("M23").Value = "TOTAL"
$ws.Range("N23").Value = 202510
$ws.Range("P23").Value = "encargado de pago"

# --- Row 24 ---
$ws.Range("A24").Value = 76477124
$ws.Range("B24").Value = "LESLY ZARATE"
$ws.Range("C24").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D24").Value = "PREJUDICIAL FLUJO"
$ws.Range("E24").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F24").Value = 20603914296
$ws.Range("H24").Value = 20603914296
$ws.Range("I24").Value = 46003
$ws.Range("L24").Value = 205.38
$ws.Range("M24").Value = "TOTAL"
$ws.Range("N24").Value = 202510

# --- Row 25 ---
$ws.Range("A25").Value = 76479084
$ws.Range("B25").Value = "LAURA VILLANUEVA"
$ws.Range("C25").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D25").Value = "PREJUDICIAL FLUJO"
$ws.Range("E25").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F25").Value = 20612463591
$ws.Range("H25").Value = 20612463591
$ws.Range("I25").Value = 46002
$ws.Range("L25").Value = 127
$ws.Range("M25").Value = "SIN PAGO"
$ws.Range("O25").Value = 202510
$ws.Range("P25").Value = "TORREJON REYES ROSA MICAELA"

# --- Row 26 ---
$ws.Range("A26").Value = 76475076
$ws.Range("B26").Value = "LAURA VILLANUEVA"
$ws.Range("C26").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D26").Value = "PREJUDICIAL FLUJO"
$ws.Range("E26").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F26").Value = 20521886871
$ws.Range("H26").Value = 20521886871
$ws.Range("I26").Value = 46003
$ws.Range("L26").Value = 174
$ws.Range("M26").Value = "TOTAL"
$ws.Range("N26").Value = 202510
$ws.Range("P26").Value = "NARVAEZ YSELA MIRIAM"

# --- Row 27 ---
$ws.Range("A27").Value = 76474019
$ws.Range("B27").Value = "LAURA VILLANUEVA"
$ws.Range("C27").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D27").Value = "PREJUDICIAL FLUJO"
$ws.Range("E27").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F27").Value = 41977629
$ws.Range("H27").Value = 10419776292
$ws.Range("I27").Value = 46003
$ws.Range("L27").Value = 102
$ws.Range("M27").Value = "PARCIAL"
$ws.Range("N27").Value = 202510
$ws.Range("O27").Value = "202507|202508|202509"
$ws.Range("P27").Value = "RABANAL MISARI ANA MARIA"

# --- Row 28 ---
$ws.Range("A28").Value = 76477631
$ws.Range("B28").Value = "CARLA CASTILLO"
$ws.Range("C28").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D28").Value = "PREJUDICIAL FLUJO"
$ws.Range("E28").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F28").Value = 20606210176
$ws.Range("H28").Value = 20606210176
$ws.Range("I28").Value = 46003
$ws.Range("L28").Value = 66.1
$ws.Range("M28").Value = "SIN PAGO"
$ws.Range("O28").Value = 202510
$ws.Range("P28").Value = "estudiocontable_palomino@hotmail.com"

# --- Row 29 (no ID_OBLIGACION / A29 stays blank) ---
$ws.Range("B29").Value = "CARLA CASTILLO"
$ws.Range("C29").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D29").Value = "REDIRECCIONAMIENTO"
$ws.Range("E29").Value = "PRIMA_AFP_REDIRECCIONAMIENTO"
$ws.Range("F29").Value = 20519499780
$ws.Range("H29").Value = 20519499780
$ws.Range("I29").Value = 46003
$ws.Range("L29").Value = 537.66

# --- Row 30 ---
$ws.Range("A30").Value = 76473991
$ws.Range("B30").Value = "SANDRA BENAVIDES"
$ws.Range("C30").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D30").Value = "PREJUDICIAL FLUJO"
$ws.Range("E30").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F30").Value = 40721603
$ws.Range("H30").Value = 10407216038
$ws.Range("I30").Value = 46004
$ws.Range("L30").Value = 92.14
$ws.Range("M30").Value = "TOTAL"
$ws.Range("N30").Value = 202510
$ws.Range("P30").Value = "encargado de pago"

# --- Row 31 ---
$ws.Range("A31").Value = 76474057
$ws.Range("B31").Value = "CHERRY MATSON"
$ws.Range("C31").Value = "DICIEMBRE WORLD 2025"
$ws.Range("D31").Value = "PREJUDICIAL FLUJO"
$ws.Range("E31").Value = "PRIMA_AFP_PREJUDICIAL_FLUJO"
$ws.Range("F31").Value = 43371623
$ws.Range("H31").Value = 10433716235
$ws.Range("I31").Value = 46004
$ws.Range("L31").Value = 66.1
$ws.Range("M31").Value = "PARCIAL"
$ws.Range("N31").Value = 202510
$ws.Range("O31").Value = 202509

# --- Update the active selection to match the saved view state ---
$ws.Range("L31").Select()
